$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 145
$ws1.Range("F3").Value = 1331
$ws1.Range("F4").Value = 1136
$ws1.Range("F5").Value = 1026
$ws1.Range("F6").Value = 1804
$ws1.Range("F7").Value = 572
$ws1.Range("F8").Value = 1208
$ws1.Range("F9").Value = 61
$ws1.Range("F11").Value = 127
$ws1.Range("F13").Value = 73
$ws1.Range("F15").Value = 699
$ws1.Range("F16").Value = 178
$ws1.Range("F17").Value = 105
$ws1.Range("F20").Value = 331
$ws1.Range("F21").Value = 159
$ws1.Range("F22").Value = 674
$ws1.Range("F23").Value = 42
$ws1.Range("F24").Value = 650
$ws1.Range("F25").Value = 160
$ws1.Range("F27").Value = 879
$ws1.Range("F28").Value = 318
$ws1.Range("F31").Value = 275
$ws1.Range("F32").Value = 12
$ws1.Range("F33").Value = 16
$ws1.Range("F34").Value = 405

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 258

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 145
$ws4.Range("F4").Value = 1331
$ws4.Range("F5").Value = 1136
$ws4.Range("F6").Value = 1026
$ws4.Range("F7").Value = 1804
$ws4.Range("F8").Value = 572
$ws4.Range("F9").Value = 1208
$ws4.Range("F10").Value = 61
$ws4.Range("F13").Value = 127
$ws4.Range("F15").Value = 73
$ws4.Range("F17").Value = 699
$ws4.Range("F18").Value = 178
$ws4.Range("F19").Value = 105
$ws4.Range("F25").Value = 331
$ws4.Range("F27").Value = 258
$ws4.Range("F28").Value = 258
$ws4.Range("F29").Value = 159
$ws4.Range("F30").Value = 674
$ws4.Range("F31").Value = 42
$ws4.Range("F32").Value = 650
$ws4.Range("F33").Value = 160
$ws4.Range("F35").Value = 879
$ws4.Range("F36").Value = 318
$ws4.Range("F41").Value = 275
$ws4.Range("F45").Value = 12
$ws4.Range("F46").Value = 16
$ws4.Range("F48").Value = 405
